# Test data added for Italy
#
# - Select cell B4 on the "Slovakia" sheet (it was previously on C11).
# - Add a new "Italy" sheet, cloned from "Belgium" (same A1:D10 layout,
#   column widths, merged cells and styles), placed after "Slovakia".
# - Fill in the Italy-specific values and make it the active sheet with
#   B4 selected.

$wb = $excel.ActiveWorkbook

# Move the selection on Slovakia to B4 before switching sheets, matching
# the recorded selection state for that sheet in the edited workbook.
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Activate()
$slovakia.Range("B4").Select() | Out-Null

# Clone the Belgium sheet (same shape as the new Italy sheet) and place it
# as the last tab, right after Slovakia.
$template = $wb.Worksheets.Item("Belgium")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Type]::Missing, $lastSheet) | Out-Null

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Fill B4 before B2 so new shared-string entries land in the same order
# as the target workbook (NGC-3145/T2219 then Italy Market).
$italy.Range("B4").Value = "NGC-3145/T2219"
$italy.Range("B2").Value = "Italy Market"

# Italy becomes the active sheet with B4 selected.
$italy.Range("B4").Select() | Out-Null
